# io_mux_cfg.xlsx update: expand iomux_cfg0 / iomux_cfg1 tables to 3-bit
# select codes (000-111) and fill in the previously unused func/dir/default_i
# values, add a "gpio[5]" / "uart_rx" row that used to be a placeholder "-".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the existing table rows completely so no stale cells / row spans are
# left behind, then rebuild them from scratch with the new layout.
$ws.Rows("2:12").Delete()

# ---- iomux_cfg0 table -------------------------------------------------
$ws.Range("A2").Value = "iomux_cfg0"

$ws.Range("A3").Value = "bits"
$ws.Range("B3").Value = "io_name"
$ws.Range("C3").Value = "iomux_cfg"
$ws.Range("D3").Value = "func"
$ws.Range("E3").Value = "dir"
$ws.Range("F3").Value = "default_i"

$ws.Range("A4").Value = "6:4"
$ws.Range("B4").Value = "SPI_MISO"
$ws.Range("C4").Value = "000"
$ws.Range("D4").Value = "spi_miso"
$ws.Range("E4").Value = "C/X"
$ws.Range("F4").Value = "0"

$ws.Range("C5").Value = "001"
$ws.Range("D5").Value = "uart_tx"
$ws.Range("E5").Value = "X/C"
$ws.Range("F5").Value = "-"

$ws.Range("C6").Value = "010"
$ws.Range("D6").Value = "gpio[5]"
$ws.Range("E6").Value = "C/C"
$ws.Range("F6").Value = "0"

$ws.Range("C7").Value = "011"
$ws.Range("D7").Value = "uart_rx"
$ws.Range("E7").Value = "C/X"
$ws.Range("F7").Value = "1"

$ws.Range("C8").Value = "100"
$ws.Range("D8").Value = "-"

$ws.Range("C9").Value = "101"
$ws.Range("D9").Value = "-"

$ws.Range("C10").Value = "110"
$ws.Range("D10").Value = "-"

$ws.Range("C11").Value = "111"
$ws.Range("D11").Value = "-"

# ---- iomux_cfg1 table -------------------------------------------------
$ws.Range("A13").Value = "iomux_cfg1"

$ws.Range("A14").Value = "bits"
$ws.Range("B14").Value = "io_name"
$ws.Range("C14").Value = "iomux_cfg"
$ws.Range("D14").Value = "func"
$ws.Range("E14").Value = "dir"
$ws.Range("F14").Value = "default_i"

$ws.Range("A15").Value = "3:1"
$ws.Range("B15").Value = "SPI_MOSI"
$ws.Range("C15").Value = "000"
$ws.Range("D15").Value = "spi_mosi"
$ws.Range("E15").Value = "X/C"
$ws.Range("F15").Value = "-"

$ws.Range("C16").Value = "001"
$ws.Range("D16").Value = "gpio[6]"
$ws.Range("E16").Value = "C/C"
$ws.Range("F16").Value = "1"

$ws.Range("C17").Value = "010"
$ws.Range("D17").Value = "uart_tx"
$ws.Range("E17").Value = "X/C"
$ws.Range("F17").Value = "-"

$ws.Range("C18").Value = "011"
$ws.Range("D18").Value = "uart_rx"
$ws.Range("E18").Value = "C/X"
$ws.Range("F18").Value = "1"

$ws.Range("C19").Value = "100"
$ws.Range("D19").Value = "-"

$ws.Range("C20").Value = "101"
$ws.Range("D20").Value = "-"

$ws.Range("C21").Value = "110"
$ws.Range("D21").Value = "-"

$ws.Range("C22").Value = "111"
$ws.Range("D22").Value = "-"

[void]$ws.Range("F19").Select()
